$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row for "005366255 / RAPHAELA / 2406.88" right before the
#    existing "004369172 / LUIZA" row (row 6).
$ws.Rows.Item(6).Insert()

# Column A must stay text (keep the leading zeros), so force a text format
# before assigning, then clear the formatting again so the cell ends up
# looking like its untouched siblings (no explicit style applied).
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "005366255"
$ws.Range("A6").ClearFormats()

$ws.Range("B6").Value = "RAPHAELA"
$ws.Range("C6").Value = 2406.88

# 2. Remove the old "005366255 / RAPHAELA / 406.88" row. After the insert
#    above, it shifted down one row, from 49 to 50.
$ws.Rows.Item(50).Delete()

# 3. Remove the "004487016 / ROGERIO / 117.22" row entirely. It originally
#    sat at row 133; the insert shifted it to 134, and the delete above
#    shifted it back up to 133.
$ws.Rows.Item(133).Delete()
